# Update the "想去人数" (F column) figures on the "展览" and "全部类型"
# sheets to the refreshed crawl counts from the gh-pages data regeneration.

$wb = $excel.ActiveWorkbook

$newValues = @{
    3  = 1434
    5  = 263
    9  = 197
    12 = 4768
    14 = 7032
    20 = 4
    21 = 4183
    22 = 1152
    23 = 82
    25 = 2766
    27 = 555
    29 = 390
    30 = 392
    31 = 411
    32 = 248
    33 = 52
    34 = 1652
    35 = 1063
    36 = 73
    37 = 760
    40 = 11
    41 = 501
    43 = 22
    44 = 93
    45 = 304
    47 = 24
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $newValues.Keys) {
        $ws.Range("F$row").Value = $newValues[$row]
    }
}
